$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift existing rows down to make room for new header row 7 ---
# Move old row 20 (A:B) down to row 21
$ws.Range("A20:B20").Cut($ws.Range("A21:B21"))

# Move old rows 15..7 down to 16..8 (bottom-up so we don't clobber data).
# Each row is cut using its exact last-used column so we don't materialize
# empty cells across the whole A:BJ width.
$lastCols = @{7="U"; 8="Q"; 9="M"; 10="E"; 11="W"; 12="H"; 13="AC"; 14="V"; 15="G"}
for ($r = 15; $r -ge 7; $r--) {
    $src = $r
    $dst = $r + 1
    $col = $lastCols[$r]
    $ws.Range("A$src`:$col$src").Cut($ws.Range("A$dst`:$col$dst"))
}

# Clean up leftover cells from the shift (source row 7 was never overwritten,
# and row 20's old cells remain after its content moved to row 21)
$ws.Range("A7:BJ7").Clear()
$ws.Range("A20:B20").Clear()

# --- Populate new row 7 with the new "cin7_aged_trim" header block ---
# Set the non-A cells first, then A7, so new shared-string indices are
# allocated in the same order Excel produced them (Barcode, Label, Box, cin7_aged_trim).
$ws.Range("B7").Value2 = "Barcode"
$ws.Range("C7").Value2 = "Label"
$ws.Range("D7").Value2 = "Box"
$ws.Range("A7").Value2 = "cin7_aged_trim"

# Apply the column-A header style (matches the style used by the other rows in col A)
$ws.Range("A8").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New blank styled row 17 (kept as a spacer, matching style used in column A) ---
$ws.Range("A8").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to match the new active cell
$ws.Range("A7").Select()

Write-Host "Dimension after: $($ws.UsedRange.Address())"
